# Rework the nomination-summary sheet:
#  - Row 34 ("Summary" header, no value) is removed entirely, which shifts the
#    five "Total ..." rows below it up by one row (their values/number formats
#    move with them and already line up with the new layout).
#  - Every remaining row's label (column A) is rewritten to the new wording
#    (the "Civilian, ...", "Air Force, ...", etc. style used throughout the
#    new sheet), while the numeric values in column B are left untouched.
#  - The new row 34 ("Total new nominations") needs a value that previously
#    didn't exist for that row; after the deletion/shift it holds the old
#    "carried over" value (0), which is exactly what the new layout expects.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old "Summary" divider row - everything below shifts up one row.
$ws.Rows(34).Delete()

# New labels for column A, row by row.
$labels = @{
    1  = "Labels"
    2  = "Congress"
    3  = "Session"
    4  = "Start Date"
    5  = "End Date"
    6  = "Civilian"
    7  = "     Civilian, New nominations"
    8  = "     Civilian, Confirmed "
    9  = "     Civilian, Unconfirmed "
    10 = "     Civilian, Withdrawn "
    11 = "     Civilian, Returned to White House "
    12 = "Other Civilian"
    13 = "     Other Civilian, New nominations"
    14 = "     Other Civilian, Confirmed "
    15 = "     Other Civilian, Unconfirmed "
    16 = "Air Force"
    17 = "     Air Force, New nominations"
    18 = "     Air Force, Confirmed "
    19 = "     Air Force, Unconfirmed "
    20 = "     Air Force, Returned to White House "
    21 = "Army"
    22 = "     Army, New nominations"
    23 = "     Army, Confirmed "
    24 = "     Army, Unconfirmed "
    25 = "     Army, Returned to White House "
    26 = "Navy"
    27 = "     Navy, New nominations"
    28 = "     Navy, Confirmed "
    29 = "     Navy, Unconfirmed "
    30 = "Marine Corps"
    31 = "     Marine Corps, New nominations"
    32 = "     Marine Corps, Confirmed "
    33 = "     Marine Corps, Unconfirmed "
    34 = "Total new nominations"
    35 = "Total carryover nominations"
    36 = "Total confirmed "
    37 = "Total unconfirmed "
    38 = "Total withdrawn "
    39 = "Total returned to the White House "
}

foreach ($row in $labels.Keys) {
    $ws.Cells.Item($row, 1).Value = $labels[$row]
}

# The new "Total new nominations" row needs a numeric value in column B;
# after the row-34 deletion this cell already holds the old carry-over
# value (0) with the right number format, but set it explicitly to be safe.
$ws.Cells.Item(34, 2).Value = 0
